# implement income page & expense page: download utility for transaction
# overview section. Expand the expense table with the new transactions and
# re-point "Rent" to the bottom of the list with its updated amount/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing date formatting (style) used by C2 so the new date
# cells keep the same built-in date number format instead of Excel minting
# a brand new custom numFmt entry.
$ws.Range("C2").Copy() | Out-Null

# New data rows (row 1 headers stay untouched).
$data = @(
    @("Spotify + Netflix Subscription", 30,  46010.29180555556),
    @("College Tuition",                500, 46005.29180555556),
    @("Food",                           200, 46000.29180555556),
    @("Rent",                           150, 45997.29180555556)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]

    $dateCell = $ws.Cells.Item($row, 3)
    $dateCell.PasteSpecial(-4122) | Out-Null
    $dateCell.Value = $item[2]
}

$excel.CutCopyMode = $false
